$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old sample/test row 2 data that is being replaced
$ws.Range("A2").ClearContents()

# Row 2
$ws.Range("B2").Value = "JRS6561"
$ws.Range("C2").Value = "SEAT"
$ws.Range("D2").Value = "2024-09-10 07:19 AM"
$ws.Range("E2").Value = "2024-09-12 23:00:00"

# Row 3 (overwrite previous row 3 sample data)
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "PYR831E"
$ws.Range("C3").Value = "TOYOTA"
$ws.Range("D3").Value = "2024-09-10 08:56 AM"
$ws.Range("E3").Value = "2024-09-13 00:00:00"

# Row 4
$ws.Range("B4").Value = "23N050"
$ws.Range("C4").Value = "MAZDA"
$ws.Range("D4").Value = "2024-09-10 09:10 AM"
$ws.Range("E4").Value = "2024-09-13 13:00:00"

# Row 5
$ws.Range("B5").Value = "NAP068A"
$ws.Range("C5").Value = "NISSAN"
$ws.Range("D5").Value = "2024-09-10 08:02 AM"
$ws.Range("E5").Value = "2024-09-14 10:00:00"

# Row 6
$ws.Range("B6").Value = "JSL2080"
$ws.Range("C6").Value = "MAZDA"
$ws.Range("D6").Value = "2024-09-10 09:18 AM"
$ws.Range("E6").Value = "2024-09-15 13:00:00"

# Row 7
$ws.Range("B7").Value = "C22BHG"
$ws.Range("C7").Value = "AUDI"
$ws.Range("D7").Value = "2024-09-10 08:53 AM"
$ws.Range("E7").Value = "2024-09-15 20:00:00"

# Row 8 (no salida value for this one)
$ws.Range("A8").Value = "ADOLFO REYES AGUIRRE"
$ws.Range("B8").Value = "JSB4919"
$ws.Range("C8").Value = "honda"
$ws.Range("D8").Value = "2024-09-11 17:41 PM"

# Apply the autofilter over the full data range
$ws.Range("A1:E8").AutoFilter()

# Register the hidden _FilterDatabase defined name that Excel creates
# for the worksheet's autofilter range
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Sheet1'!`$A`$1:`$E`$8")
$filterName.Visible = $false
